$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new header columns O1:Q1 ---
$ws.Range("O1").Value = "coverage_error"
$ws.Range("P1").Value = "label_ranking_average_precision_score"
$ws.Range("Q1").Value = "label_ranking_loss"

# Match the header formatting (bold/centered/bordered) used by the rest of row 1
$ws.Range("N1").Copy()
$ws.Range("O1:Q1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# --- Extend the existing data rows (2-10) with blank cells in the new columns ---
$ws.Range("O2:Q10").Style = "Normal"

# --- Add new row 11 for the "tpot" framework ---
# Columns B..L and O..Q hold numeric-looking metric strings (e.g. "0.983",
# "-1.000") that must stay TEXT (to keep the trailing zeros / exact digits,
# matching every other data row). Force text entry via a temporary "@"
# number format, then drop back to the default style so no visible
# formatting change is introduced.
$textCells = "B11","C11","D11","E11","F11","G11","H11","I11","J11","K11","L11","O11","P11","Q11"
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("A11").Value = "tpot"
$ws.Range("B11").Value = "0.983"
$ws.Range("C11").Value = "-1.000"
$ws.Range("D11").Value = "0.947"
$ws.Range("E11").Value = "0.964"
$ws.Range("F11").Value = "0.923"
$ws.Range("G11").Value = "0.983"
$ws.Range("H11").Value = "0.984"
$ws.Range("I11").Value = "0.964"
$ws.Range("J11").Value = "-1.000"
$ws.Range("K11").Value = "-1.000"
$ws.Range("L11").Value = "-1.000"
$ws.Range("M11").Value = "00:10:10"
$ws.Range("N11").Value = "00:00:00"
$ws.Range("O11").Value = "-1.000"
$ws.Range("P11").Value = "-1.000"
$ws.Range("Q11").Value = "-1.000"

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}

$wb.Save()
